# Regenerate merged AHB files
# Rename the "_old"/"_new" header suffixes to "_FV2404"/"_FV2410", turn the
# merged-AHB data range into an Excel Table (Table1) and freeze the header
# row, matching the regenerated output of the AHB-merge tooling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

# --- 1. Rename header row labels (row 1) ----------------------------------
# "<name>_old" -> "<name>_FV2404"   (previous/base AHB version columns)
# "<name>_new" -> "<name>_FV2410"   (current/target AHB version columns)
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2404")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2410")
    }
}

# --- 2. Create a table (ListObject) over the data with autofilter --------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$listObject = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$listObject.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true
